$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:B1").EntireColumn.Insert()

$rng = $ws1.Range("A2:A1048576")
$rng.Validation.Add(3, 1, 1, "'version list'!`$A`$1:`$A`$1")
Write-Output "added"
